$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Teste2: turn the single "Acura" cell into a full table like Teste1 ---
$ws2.Range("A1").Value = "firstname"
$ws2.Range("B1").Value = "lastname"
$ws2.Range("C1").Value = "usarname"
$ws2.Range("D1").Value = "password"

$ws2.Range("A2").Value = "Francisco"
$ws2.Range("B2").Value = "Zimmer"
$ws2.Range("C2").Value = "zimmerf "
$ws2.Range("D2").Value = 852456

$ws2.Range("A3").Value = "Felipe"
$ws2.Range("B3").Value = "Zimmer"
$ws2.Range("C3").Value = "zimmerf "
$ws2.Range("D3").Value = 123456

# --- Selection / active-sheet bookkeeping ---
# Teste1 ends up with A1:D4 selected (no longer the active tab).
$ws1.Range("A1:D4").Select()

# Teste2 becomes the active tab, with G5 selected.
$ws2.Activate()
$ws2.Range("G5").Select()
